# Regenerate save_data column G ("K") values: use K instead of Strike#,
# i.e. overwrite the previously-stored Strike# values with the newly
# computed K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 2
    8  = 0
    9  = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 3
    15 = 1
    17 = 2
    18 = 2
    19 = 1
    20 = 0
    21 = 0
    22 = 1
    23 = 3
    24 = 0
    25 = 2
    26 = 1
    27 = 0
    28 = 1
    29 = 2
    30 = 1
    31 = 2
    32 = 2
    33 = 2
    34 = 1
    35 = 0
    36 = 0
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 2
    45 = 1
    47 = 1
    48 = 1
    49 = 1
    50 = 1
    51 = 0
    52 = 1
    53 = 1
    54 = 1
    55 = 1
    57 = 1
    58 = 1
    59 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
